$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Certifications")
$ws.Range("A78:G78").Copy()
$ws.Range("A79:G79").PasteSpecial(-4123)  # xlPasteAll? let's try
Write-Host "pasted all"
